# "Generate Report for Archive" - refresh the localization-status report:
#  - the single in-flight item moved from "Ready for handoff" to "In Translation"
#  - the Overview zh-cn/de-de columns and the per-locale Status columns were
#    re-sized (narrower) by the report generator to fit the new status text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: E2 (zh-cn status) and F2 (de-de status) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Columns E and F narrow to fit the shorter status text.
$overview.Columns.Item(5).ColumnWidth = 12.57
$overview.Columns.Item(6).ColumnWidth = 12.57

# --- zh-cn sheet: Status column (C2) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.57

# --- de-de sheet: Status column (C2) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.57
